# Update cryptos list — 2024-03-27 refresh (values + a couple of row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number need an explicit
# Text number format first, otherwise Excel auto-converts the assigned
# string into a numeric value (dropping e.g. a trailing '.0').
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.639.76'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '3.483.63'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '569.10'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").Value = '182.32'
$ws.Range("E6").Value = '  -3.78%  '
$ws.Range("E7").Value = '  -3.33%  '
$ws.Range("D8").Value = '3.478.89'
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '0.183'
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  -3.86%  '
$ws.Range("D12").Value = '53.55'
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = '4.044.82'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '19.14'
$ws.Range("E16").Value = '  -3.88%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.489.80'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '68.597.63'
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("D19").Value = '12.27'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("E20").Value = '  -1.60%  '
$ws.Range("D21").Value = '538.65'
$ws.Range("E21").Value = '  +12.69%  '
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("D23").Value = '19.40'
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("D24").Value = '4.96'
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").Value = '93.86'
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("D27").Value = '2.89'
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("D28").Value = '10.77'
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("D30").Value = '31.26'
$ws.Range("E30").Value = '  -3.19%  '
$ws.Range("E31").Value = '  -6.96%  '
$ws.Range("D32").Value = '12.54'
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").Value = '64.08'
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("E34").Value = '  -5.49%  '
$ws.Range("D35").Value = '568.67'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '37.66'
$ws.Range("E37").Value = '  -3.77%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '3.01'
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = '0.394'
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  -4.86%  '
$ws.Range("D41").Value = '0.133'
$ws.Range("E41").Value = '  -4.92%  '
$ws.Range("D42").Value = '3.06'
$ws.Range("E42").Value = '  -6.80%  '
$ws.Range("E43").Value = '  -5.33%  '
$ws.Range("D44").Value = '3.199.85'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("E45").Value = '  -4.84%  '
$ws.Range("D46").Value = '3.45'
$ws.Range("E46").Value = '  +3.20%  '
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("D48").Value = '9.02'
$ws.Range("E48").Value = '  -4.61%  '
$ws.Range("D49").Value = '0.134'
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("D51").Value = '136.19'
$ws.Range("E51").Value = '  -0.56%  '
